$wb = $excel.ActiveWorkbook

$wsPsgr = $wb.Worksheets.Item("SoCDTtiNTY-psgr")
$wsFrgt = $wb.Worksheets.Item("SoCDTtiNTY-frgt")

# --- SoCDTtiNTY-psgr: replace calibration formulas with hard-coded values ---
$wsPsgr.Range("B2:H2").Value = 0.086
$wsPsgr.Range("B3:H3").Value = 0.09
$wsPsgr.Range("B4:H4").Value = 0.0416
$wsPsgr.Range("B5:H5").Value = 0.029
$wsPsgr.Range("B6:H6").Value = 0.02982
$wsPsgr.Range("B7:H7").Value = 0.0587
$wsPsgr.Range("D7").Value = 0.068

# --- SoCDTtiNTY-frgt: replace calibration formulas with hard-coded values ---
$wsFrgt.Range("B2:H2").Value = 0.072
$wsFrgt.Range("B3:H3").Value = 0.0355
$wsFrgt.Range("B4:H4").Value = 0.028
$wsFrgt.Range("B5:H5").Value = 0.029
$wsFrgt.Range("B6:H6").Value = 0.0303

# --- Restore the active cell selections left on each sheet, then leave
#     "About" as the active tab (matches the saved state in the workbook) ---
$wsPsgr.Range("B7:H7").Select() | Out-Null
$wsFrgt.Range("B4:H4").Select() | Out-Null

$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate() | Out-Null
